# Updated cryptos list data (Price + Volume(1h) columns) per latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.318.25"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.03%  "

$ws.Range("D3").Value = "'1.928.57"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.22%  "

$ws.Range("D4").Value = "'1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.11%  "

$ws.Range("D5").Value = "'0.7500"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.21%  "

$ws.Range("D6").Value = "'243.80"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.96%  "

$ws.Range("E7").Value = "  +0.18%  "

$ws.Range("D8").Value = "'0.3153"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.27%  "

$ws.Range("D9").Value = "'27.39"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.19%  "

$ws.Range("D10").Value = "'0.06965"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.18%  "

$ws.Range("D11").Value = "'0.08003"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.14%  "

$ws.Range("D12").Value = "'0.7703"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.52%  "

$ws.Range("D13").Value = "'1.932.15"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.00%  "

$ws.Range("D14").Value = "'5.318"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.16%  "

$ws.Range("D15").Value = "'93.21"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.71%  "

$ws.Range("D16").Value = "'30.317.11"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.04%  "

$ws.Range("D17").Value = "'14.30"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.53%  "

$ws.Range("D18").Value = "'250.78"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.83%  "

$ws.Range("D19").Value = "'0.000007880"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.78%  "

$ws.Range("D20").Value = "'5.734"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.09%  "

$ws.Range("D21").Value = "'2.187.92"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.38%  "

$ws.Range("D23").Value = "'1.002"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.06%  "

$ws.Range("D24").Value = "'6.629"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.19%  "

$ws.Range("D25").Value = "'9.413"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.98%  "

$ws.Range("D26").Value = "'165.71"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.56%  "

$ws.Range("D27").Value = "'18.89"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.95%  "

$ws.Range("D28").Value = "'0.1326"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.50%  "

$ws.Range("D29").Value = "'2.185"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.59%  "

$ws.Range("D30").Value = "'1.370"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.65%  "

$ws.Range("E31").Value = "  -2.04%  "

$ws.Range("D32").Value = "'4.366"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.52%  "

$ws.Range("D33").Value = "'4.086"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.73%  "

$ws.Range("D34").Value = "'0.05105"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.96%  "

$ws.Range("D35").Value = "'1.276"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.82%  "

$ws.Range("D36").Value = "'0.7426"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.94%  "

$ws.Range("D37").Value = "'2.773"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.20%  "

$ws.Range("D38").Value = "'0.01944"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.92%  "

$ws.Range("D39").Value = "'2.793"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.33%  "

$ws.Range("D40").Value = "'76.83"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.79%  "

$ws.Range("D41").Value = "'6.363"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.09%  "

$ws.Range("D42").Value = "'0.4420"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.37%  "

$ws.Range("D43").Value = "'1.949"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.78%  "

$ws.Range("D44").Value = "'1.002"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.08%  "

$ws.Range("D45").Value = "'0.8304"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.41%  "

$ws.Range("D46").Value = "'100.45"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.39%  "

$ws.Range("D47").Value = "'9.714"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.92%  "

$ws.Range("D48").Value = "'7.426"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.27%  "

$ws.Range("D49").Value = "'37.14"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.19%  "

$ws.Range("D50").Value = "'973.06"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +9.34%  "

$ws.Range("D51").Value = "'0.06041"
$ws.Range("D51").Style = "Normal"
